# "Update Day 2 and 3"
#
# Adds a new worksheet "S2" right after "TSLA_EN" holding a snapshot of the
# Date/Open/High/Low/Close columns (A:E) for the header row plus the first
# 19 trading days (rows 1-20) of TSLA_EN. The new sheet becomes the active /
# visible tab, and both sheets end up with the A1:E20 block selected.

$wb = $excel.ActiveWorkbook
$srcName = "TSLA_EN"
$newName = "S2"

$srcForAdd = $wb.Worksheets.Item($srcName)

# Insert the new sheet immediately after TSLA_EN.
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $srcForAdd)
$new.Name = $newName

# NOTE: worksheet object handles obtained before Worksheets.Add() go stale,
# so re-resolve both sheets by name before touching them again.
$src = $wb.Worksheets.Item($srcName)
$dst = $wb.Worksheets.Item($newName)

# Match number formats first (column data is stored as text, e.g. "2018-12-12"
# and "369.420013" -- copying the cell format before the value keeps Excel
# from auto-converting those text-looking-like-dates/numbers on write).
$dst.Range("A1:E20").NumberFormat = $src.Range("A1:E20").NumberFormat

for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $dst.Cells.Item($r, $c).Value2 = $src.Cells.Item($r, $c).Value2
    }
}

# Keep page setup consistent with the source sheet.
$dst.PageSetup.LeftMargin = $src.PageSetup.LeftMargin
$dst.PageSetup.RightMargin = $src.PageSetup.RightMargin
$dst.PageSetup.TopMargin = $src.PageSetup.TopMargin
$dst.PageSetup.BottomMargin = $src.PageSetup.BottomMargin
$dst.PageSetup.HeaderMargin = $src.PageSetup.HeaderMargin
$dst.PageSetup.FooterMargin = $src.PageSetup.FooterMargin

# Leave the same A1:E20 block selected on both sheets ...
[void]$src.Range("A1:E20").Select()
[void]$dst.Range("A1:E20").Select()

# ... with the new sheet as the active / visible tab.
[void]$dst.Activate()
